$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 14686966
$ws.Range("I19").Value = 18184272
$ws.Range("J19").Value = 12501150
$ws.Range("K19").Value = 18184272
$ws.Range("L19").Value = 12501150
$ws.Range("M19").Value = -18184097
$ws.Range("N19").Value = -12501500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1605.2916
$ws.Range("I2").Value = 1396.3889
$ws.Range("J2").Value = 2232
$ws.Range("K2").Value = 1396.3889
$ws.Range("L2").Value = 2232
$ws.Range("M2").Value = -1283.3889
$ws.Range("N2").Value = -2458
$ws.Range("H45").Value = 1174.5454
$ws.Range("I45").Value = 565.3333
$ws.Range("J45").Value = 1905.6
$ws.Range("K45").Value = 565.3333
$ws.Range("L45").Value = 1905.6
$ws.Range("M45").Value = -188.3333
$ws.Range("N45").Value = -2659.6
$ws.Range("H74").Value = 1645.4546
$ws.Range("I74").Value = 1417.037
$ws.Range("J74").Value = 2673.3333
$ws.Range("K74").Value = 1417.037
$ws.Range("L74").Value = 2673.3333
$ws.Range("M74").Value = -543.037
$ws.Range("N74").Value = -4421.3333
$ws.Range("H77").Value = 1645.4546
$ws.Range("I77").Value = 1417.037
$ws.Range("J77").Value = 2673.3333
$ws.Range("K77").Value = 7085.185
$ws.Range("L77").Value = 13366.6665
$ws.Range("M77").Value = -2717.185
$ws.Range("N77").Value = -22102.6665
$ws.Range("H116").Value = 1605.2916
$ws.Range("I116").Value = 1396.3889
$ws.Range("J116").Value = 2232
$ws.Range("K116").Value = 1396.3889
$ws.Range("L116").Value = 2232
$ws.Range("M116").Value = 897.6111000000001
$ws.Range("N116").Value = -6820
$ws.Range("H132").Value = 2294.7827
$ws.Range("I132").Value = 2442.3076
$ws.Range("J132").Value = 2103
$ws.Range("K132").Value = 7326.9228
$ws.Range("L132").Value = 6309
$ws.Range("M132").Value = -4796.9228
$ws.Range("N132").Value = -11369

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1605.2916
$ws.Range("I3").Value = 1396.3889
$ws.Range("J3").Value = 2232
$ws.Range("K3").Value = 1396.3889
$ws.Range("L3").Value = 2232
$ws.Range("M3").Value = -1282.3889
$ws.Range("N3").Value = -2460
$ws.Range("H80").Value = 182.88889
$ws.Range("I80").Value = 175.8
$ws.Range("J80").Value = 191.75
$ws.Range("K80").Value = 175.8
$ws.Range("L80").Value = 191.75
$ws.Range("M80").Value = 822.2
$ws.Range("N80").Value = -2187.75
$ws.Range("H83").Value = 182.88889
$ws.Range("I83").Value = 175.8
$ws.Range("J83").Value = 191.75
$ws.Range("K83").Value = 879
$ws.Range("L83").Value = 958.75
$ws.Range("M83").Value = 4113
$ws.Range("N83").Value = -10942.75
$ws.Range("H94").Value = 385.4
$ws.Range("I94").Value = 339.33334
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 339.33334
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 111.66666
$ws.Range("N94").Value = -1702
$ws.Range("H99").Value = 7101.579
$ws.Range("I99").Value = 11777.556
$ws.Range("J99").Value = 2893.2
$ws.Range("K99").Value = 11777.556
$ws.Range("L99").Value = 2893.2
$ws.Range("M99").Value = -10279.556
$ws.Range("N99").Value = -5889.2
$ws.Range("H134").Value = 2036
$ws.Range("I134").Value = 1838.8572
$ws.Range("J134").Value = 2537.818
$ws.Range("K134").Value = 5516.571599999999
$ws.Range("L134").Value = 7613.454000000001
$ws.Range("M134").Value = -2981.571599999999
$ws.Range("N134").Value = -12683.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3881.8965
$ws.Range("I31").Value = 2358.4375
$ws.Range("J31").Value = 5756.923
$ws.Range("K31").Value = 2358.4375
$ws.Range("L31").Value = 5756.923
$ws.Range("M31").Value = -2063.4375
$ws.Range("N31").Value = -6346.923
$ws.Range("H34").Value = 3881.8965
$ws.Range("I34").Value = 2358.4375
$ws.Range("J34").Value = 5756.923
$ws.Range("K34").Value = 2358.4375
$ws.Range("L34").Value = 5756.923
$ws.Range("M34").Value = -2156.4375
$ws.Range("N34").Value = -6160.923
$ws.Range("H134").Value = 2687.4285
$ws.Range("I134").Value = 1644.05
$ws.Range("J134").Value = 5295.875
$ws.Range("K134").Value = 4932.15
$ws.Range("L134").Value = 15887.625
$ws.Range("M134").Value = -2397.15
$ws.Range("N134").Value = -20957.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1222.7778
$ws.Range("I122").Value = 1381.4
$ws.Range("J122").Value = 1024.5
$ws.Range("K122").Value = 4144.200000000001
$ws.Range("L122").Value = 3073.5
$ws.Range("M122").Value = -1694.200000000001
$ws.Range("N122").Value = -7973.5
$ws.Range("H126").Value = 3573237.5
$ws.Range("I126").Value = 5556774
$ws.Range("J126").Value = 2872.2
$ws.Range("K126").Value = 16670322
$ws.Range("L126").Value = 8616.599999999999
$ws.Range("M126").Value = -16667852
$ws.Range("N126").Value = -13556.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4690.517
$ws.Range("I40").Value = 3657.6667
$ws.Range("J40").Value = 5797.143
$ws.Range("K40").Value = 3657.6667
$ws.Range("L40").Value = 5797.143
$ws.Range("M40").Value = -3521.6667
$ws.Range("N40").Value = -6069.143
$ws.Range("H61").Value = 1267.6364
$ws.Range("I61").Value = 1092.3334
$ws.Range("K61").Value = 1092.3334
$ws.Range("M61").Value = -890.3334
$ws.Range("H113").Value = 1267.6364
$ws.Range("I113").Value = 1092.3334
$ws.Range("K113").Value = 1092.3334
$ws.Range("M113").Value = 1077.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2029.3513
$ws.Range("I122").Value = 1767.9231
$ws.Range("K122").Value = 5303.7693
$ws.Range("M122").Value = -2853.7693
$ws.Range("H126").Value = 3012.3333
$ws.Range("I126").Value = 2372.8572
$ws.Range("K126").Value = 7118.571599999999
$ws.Range("M126").Value = -4648.571599999999
$ws.Range("H132").Value = 1842.5264
$ws.Range("I132").Value = 1320.8
$ws.Range("J132").Value = 2182.7827
$ws.Range("K132").Value = 3962.4
$ws.Range("L132").Value = 6548.348100000001
$ws.Range("M132").Value = -1432.4
$ws.Range("N132").Value = -11608.3481
$ws.Range("H136").Value = 4432.486
$ws.Range("I136").Value = 4979.0415
$ws.Range("J136").Value = 3240
$ws.Range("K136").Value = 14937.1245
$ws.Range("L136").Value = 9720
$ws.Range("M136").Value = -12387.1245
$ws.Range("N136").Value = -14820
